# Update the dSF column (F) values for the rows that changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 0
    4  = -2
    6  = -3
    11 = -6
    13 = 4
    15 = 6
    16 = -2
    18 = 2
    19 = -7
    22 = 0
    28 = 3
    31 = -3
    32 = 5
    34 = 1
    37 = -3
    38 = -5
    41 = 5
    42 = -2
    49 = 5
    50 = 2
    60 = -4
    67 = -2
    68 = 0
    69 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
